$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.098.20"
$ws.Range("E2").Value = "  -4.81%  "
$ws.Range("D3").Value = "3.340.76"
$ws.Range("E3").Value = "  -5.73%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'560.83"
$ws.Range("E5").Value = "  -4.15%  "
$ws.Range("D6").Value = "'182.50"
$ws.Range("E6").Value = "  -7.55%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  -3.14%  "
$ws.Range("D9").Value = "3.329.64"
$ws.Range("E9").Value = "  -5.66%  "
$ws.Range("D10").Value = "'0.184"
$ws.Range("E10").Value = "  -9.61%  "
$ws.Range("E11").Value = "  -6.72%  "
$ws.Range("D12").Value = "'47.47"
$ws.Range("E12").Value = "  -8.66%  "
$ws.Range("E13").Value = "  -7.69%  "
$ws.Range("D14").Value = "'8.64"
$ws.Range("E14").Value = "  -6.72%  "
$ws.Range("D15").Value = "3.877.47"
$ws.Range("E15").Value = "  -5.64%  "
$ws.Range("D16").Value = "'605.96"
$ws.Range("E16").Value = "  -8.97%  "
$ws.Range("D17").Value = "66.224.91"
$ws.Range("E17").Value = "  -4.87%  "
$ws.Range("D18").Value = "'18.06"
$ws.Range("E18").Value = "  -2.65%  "
$ws.Range("D19").Value = "3.345.57"
$ws.Range("E19").Value = "  -5.92%  "
$ws.Range("E20").Value = "  -3.83%  "
$ws.Range("E21").Value = "  -7.94%  "
$ws.Range("E22").Value = "  -6.39%  "
$ws.Range("D23").Value = "'16.95"
$ws.Range("E23").Value = "  -7.24%  "
$ws.Range("D24").Value = "'5.03"
$ws.Range("E24").Value = "  -4.81%  "
$ws.Range("D25").Value = "'99.92"
$ws.Range("E25").Value = "  -5.31%  "
$ws.Range("D26").Value = "'4.06"
$ws.Range("E26").Value = "  -7.28%  "
$ws.Range("D27").Value = "'5.99"
$ws.Range("E27").Value = "  -0.65%  "
$ws.Range("D28").Value = "'2.68"
$ws.Range("E28").Value = "  -8.20%  "
$ws.Range("D29").Value = "'9.32"
$ws.Range("E29").Value = "  -8.48%  "
$ws.Range("E30").Value = "  -9.74%  "
$ws.Range("D31").Value = "'30.40"
$ws.Range("E31").Value = "  -9.26%  "
$ws.Range("E32").Value = "  -8.51%  "
$ws.Range("D33").Value = "'3.79"
$ws.Range("E33").Value = "  -14.99%  "
$ws.Range("D34").Value = "'11.05"
$ws.Range("E34").Value = "  -6.87%  "
$ws.Range("D35").Value = "3.857.05"
$ws.Range("E35").Value = "  +1.92%  "
$ws.Range("D36").Value = "'546.48"
$ws.Range("E36").Value = "  +8.56%  "
$ws.Range("E37").Value = "  -5.52%  "
$ws.Range("D38").Value = "'57.52"
$ws.Range("E38").Value = "  -7.08%  "
$ws.Range("D39").Value = "'0.999"
$ws.Range("E39").Value = "  -0.09%  "
$ws.Range("D40").Value = "'3.40"
$ws.Range("E40").Value = "  -8.79%  "
$ws.Range("D41").Value = "0.0₃0715"
$ws.Range("E41").Value = "  -12.25%  "
$ws.Range("E42").Value = "  -8.93%  "
$ws.Range("E43").Value = "  -6.75%  "
$ws.Range("E44").Value = "  -8.10%  "
$ws.Range("D45").Value = "'32.10"
$ws.Range("E45").Value = "  -7.49%  "
$ws.Range("E46").Value = "  +17.64%  "
$ws.Range("D47").Value = "'0.0413"
$ws.Range("E47").Value = "  -8.80%  "
$ws.Range("D48").Value = "'3.12"
$ws.Range("E48").Value = "  -7.93%  "
$ws.Range("D49").Value = "'2.62"
$ws.Range("E49").Value = "  -9.14%  "
$ws.Range("E50").Value = "  -4.77%  "
$ws.Range("E51").Value = "  -0.05%  "
